$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 17 (pond "Stick_and_Stone") column C previously held a shared formula
# (=B17/10000). Replace it with the literal text value "s" since the
# sediment-volume calc was moved to another repo (this file now only
# covers bathymetry / pond area).
$ws.Range("C17").ClearContents()
$ws.Range("C17").Value = "s"
